# Update crypto price (column D) and 1h volume/change (column E) values
# Values are kept as text (matching the original inlineStr cells) by forcing
# the NumberFormat to "@" (Text) before assigning, avoiding float round-off.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.503.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.289.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.61'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.62'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.64%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0957'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.47%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.338'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.74'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.699.68'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.92'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.482.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.296.78'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.29'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.15%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '305.13'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.62%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.10%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.95'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.54%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.33'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0698'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.10%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.04%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.965'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +10.58%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.18%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.62%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.64%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.82%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '126.33'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0496'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.61%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.90%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '242.58'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.28%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.77'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.77%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.55%  '

